# Add 2022-Q3 data:
#  1) Insert a new summary row for "2022-Q3" at the top of the "总计" sheet's
#     data (shifting the existing quarters down by one row).
#  2) Insert a new worksheet named "2022-Q3" right after "总计" containing the
#     per-fund holdings detail for that quarter.

$wb = $excel.ActiveWorkbook

# --- 1. "总计" overview sheet: insert new row 2 for 2022-Q3 -----------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.8100000000000001

# Copy the number-column style (border/alignment) from the row below so A2
# matches the look of the other rows in column A.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# --- 2. New "2022-Q3" worksheet, placed right after "总计" ------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# After insertion, the worksheet that used to be "2021-Q4" (still holding the
# header/row styling we want to copy) now sits at index 3.
$refSheet = $wb.Worksheets.Item(3)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Match header formatting used by the other quarter sheets.
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'010695"
$newSheet.Range("C2").Value = "华夏磐益一年定期开放混合"
$newSheet.Range("D2").Value = "'15.90"
$newSheet.Range("E2").Value = "'99.95"
$newSheet.Range("F2").Value = "'5.12"
$newSheet.Range("G2").Value = "'0.8141"
$newSheet.Range("H2").Value = 2

# Match the column-A formatting used by the other quarter sheets.
$refSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Keep "总计" as the active tab, as in the source workbook.
$summary.Activate()
